$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# Rename sheet 3 "catasis" -> "cat asis"
$ws3.Name = "cat asis"

# Update selection on "PCA order" sheet (sheet2): A1:D27 -> D1
$ws2.Range("D1").Select() | Out-Null

# Update selection on "cat asis" sheet (sheet3): E11 -> D12, and scroll so topLeftCell resets
$ws3.Range("D12").Select() | Out-Null

# Set new column widths on "cat asis" sheet
$ws3.Columns.Item(1).ColumnWidth = 11.5
$ws3.Columns.Item(2).ColumnWidth = 9.6666666666666667
$ws3.Columns.Item(112).ColumnWidth = 7.8333333333333333

# Make "cat asis" the active/selected tab (removes tabSelected from Sheet1, sets activeTab=2)
$ws3.Activate()
